$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows that flip the Whitelist (D) flag from FALSE to TRUE, and also
# receive a new Timestamp (E) value.
$rowsWithTimestamp = 2,5,7,9,11,13,32,37,39

foreach ($r in $rowsWithTimestamp) {
    $ws.Cells.Item($r, 4).Value = $true
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = 45987
    $eCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# Rows that flip the Whitelist (D) flag from FALSE to TRUE, but do NOT
# get a Timestamp value.
$rowsNoTimestamp = 6,8,33,34,36,38

foreach ($r in $rowsNoTimestamp) {
    $ws.Cells.Item($r, 4).Value = $true
}

# Row 35 already had Whitelist = TRUE; just refresh its Timestamp.
$ws.Cells.Item(35, 5).Value = 45987
